$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("C1").Value = "max_injection_rate_5"
$ws.Range("D1").Value = "max_injection_rate_10"

# Update data rows (B,C,D columns), row 4 (A=15) has B4,C4 cleared
$data = @(
    @(62, 74, 75),
    @(82, 99, 100),
    @($null, $null, 113),
    @(101, 111, 120),
    @(106, 118, 126),
    @(106, 127, 130),
    @(115, 128, 132),
    @(107, 119, 136),
    @(111, 136, 137),
    @(107, 124, 138),
    @(121, 125, 140),
    @(110, 128, 141),
    @(109, 132, 142),
    @(111, 133, 146),
    @(121, 131, 142),
    @(114, 130, 144),
    @(123, 139, 158),
    @(147, 149, 172),
    @(142, 161, 183)
)

$row = 2
foreach ($vals in $data) {
    if ($null -eq $vals[0]) {
        $ws.Cells.Item($row, 2).ClearContents()
    } else {
        $ws.Cells.Item($row, 2).Value = $vals[0]
    }
    if ($null -eq $vals[1]) {
        $ws.Cells.Item($row, 3).ClearContents()
    } else {
        $ws.Cells.Item($row, 3).Value = $vals[1]
    }
    if ($null -eq $vals[2]) {
        $ws.Cells.Item($row, 4).ClearContents()
    } else {
        $ws.Cells.Item($row, 4).Value = $vals[2]
    }
    $row++
}
